$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 66: rename Trigger text (was "Inspect Photo - Owner and Partner")
$ws.Range("B66").Value = 'Inspect Photo - Lobster Paul'

# Rows 68-70: yellow highlight (kept content, new styling + blank filler cells)
$yellow = 65535
$ws.Range("A68:C70").Interior.Color = $yellow
$ws.Range("D68:F69").Interior.Color = $yellow
$ws.Range("D70").Interior.Color = $yellow

# Row 71: red highlight (kept content, new styling + blank filler cells)
$red = 255
$ws.Range("A71:C71").Interior.Color = $red
$ws.Range("D71:F71").Interior.Color = $red

# New rows 77-84 (appended photo-inspection lines)
$ws.Range("A77").Value = 'All Parts'
$ws.Range("B77").Value = 'Inspect Photo - Halloween Party'
$ws.Range("C77").Value = 'Paul''s first Halloween with us. He and the guys really hit it off.'
$ws.Range("A78").Value = 'All Parts'
$ws.Range("B78").Value = 'Inspect Photo - Camping (Duo)'
$ws.Range("C78").Value = 'That was supposed to be a family photo, remember? Paul just couldn''t get you to stay in frame!'
$ws.Range("A79").Value = 'All Parts'
$ws.Range("B79").Value = 'Inspect Photo - Paul with Baby'
$ws.Range("C79").Value = 'Paul always wanted kids. I guess you and his nephew were the closest he could get.'
$ws.Range("A80").Value = 'All Parts'
$ws.Range("B80").Value = 'Inspect Photo - Crazy Paul'
$ws.Range("C80").Value = 'He really was one hell of a crazy guy. I don''t think he could sit for more than one photo without making some kinda face.'
$ws.Range("A81").Value = 'All Parts'
$ws.Range("B81").Value = 'Inspect Photo - Drinks'
$ws.Range("C81").Value = 'I think that''s the last time we went out together. It just got too hard after that.'
$ws.Range("A82").Value = 'All Parts'
$ws.Range("B82").Value = 'Inspect Photo - Oscar Truck'
$ws.Range("C82").Value = 'Damn, my hair really did look better back then.'
$ws.Range("A83").Value = 'All Parts'
$ws.Range("B83").Value = 'Inspect Photo - Camping (Group)'
$ws.Range("C83").Value = 'Not pictured - Paul running off to keep your ass from jumping in the river.'
$ws.Range("A84").Value = 'All Parts'
$ws.Range("B84").Value = 'Inspect Photo - Wedding'
$ws.Range("C84").Value = 'Betty and Jim''s wedding. See us there, upper left? You know, they kept asking us when we were gonna tie it. I should give ''em a call. It''s been too long. '

# New column header H1 = "Replaced" with yellow fill
$ws.Range("H1").Value = 'Replaced'
$ws.Range("H1").Interior.Color = $yellow
$ws.Range("H1").Font.Name = "Arial"
$ws.Range("H1").Font.Size = 10
$ws.Range("H1").Font.Bold = $false

# Apply the "retyped" font treatment used for freshly edited cells
$ws.Range("B66").Font.Name = "Arial"
$ws.Range("B66").Font.Size = 10
$ws.Range("A78:C84").Font.Name = "Arial"
$ws.Range("A78:C84").Font.Size = 10

